# Add the 2021 data column (column R) to the "3.9.2" indicator sheet.
# Mirrors the existing 2020 column (Q): same per-row formatting, new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column R, row by row (row 4 is the year header).
$values = @{
    4  = 2021
    5  = 1
    6  = 2.2
    7  = 1.7
    8  = "-"
    9  = 0.3
    10 = 1.1
    11 = "-"
    12 = 0.9
    13 = 0.4
    14 = 0.6
}

foreach ($row in 4..14) {
    # Clone formatting (number format / font / borders / alignment) from the
    # column immediately to the left (Q), which carries the equivalent style
    # for this row, then overwrite with the real 2021 value.
    $ws.Range("Q$row").Copy($ws.Range("R$row"))
    $ws.Range("R$row").Value = $values[$row]
}

# Move the remembered selection the same way the source workbook does
# (previously O17, now one column further right: S17).
$ws.Range("S17").Select()
